# Update M2/FX length and date columns to reflect the latest month of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => hashtable of cell updates (column letter => new value)
$updates = @{
    2  = @{ C = 331;  F = 45108 }
    3  = @{ E = 29921; F = 45108 }
    5  = @{ E = 29921; F = 45108 }
    7  = @{ E = 29891; F = 45078 }
    11 = @{ E = 29891; F = 45078 }
    27 = @{ E = 29921; F = 45108 }
    33 = @{ C = 451;  F = 45108 }
    34 = @{ C = 199;  F = 45108 }
    36 = @{ C = 499;  F = 45108 }
    40 = @{ C = 283;  F = 45108 }
    42 = @{ C = 235;  F = 45078 }
    44 = @{ C = 402;  F = 45078 }
    49 = @{ C = 294;  F = 45078 }
    50 = @{ C = 356;  F = 45108 }
    52 = @{ C = 343;  F = 45108 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
